# Update countries & provincias Spain
# - Refresh case counts for a handful of countries.
# - Re-sort the table by "Casos totales" (column B) descending, as the
#   page always keeps the ranking in sync with the latest numbers.
# - Bump the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the case counts for the countries whose figures changed ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

$updates = @(
    @{ Row = 6;  B = 38757; C = 14550; D = 178;  E = 38179; F = 708;  G = 98;  H = 400 },  # Estados Unidos
    @{ Row = 10; B = 16018; C = 1559;  D = 1587; E = 13757; F = 1746; G = 112; H = 674 },  # Francia
    @{ Row = 12; B = 7474;  C = 611;   D = 131;  E = 7245;  F = 141;  G = 18;  H = 98  },  # Suiza
    @{ Row = 31; B = 906;   C = 121;   D = 5;    E = 897;   F = 29;   G = 1;   H = 4   },  # Irlanda
    @{ Row = 33; B = 646;   C = 1;     D = 13;   E = 629;   F = 0;    G = 1;   H = 4   },  # Pakistan
    @{ Row = 58; B = 274;   C = 34;    D = 2;    E = 272;   F = 0;    G = 0;   H = 0   },  # Sudafrica
    @{ Row = 82; B = 112;   C = 12;    D = 1;    E = 111;   F = 0;    G = 0;   H = 0   }   # Jordania
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# --- 2) Re-sort the data rows (A4:H192) by "Casos totales" descending ---
$sortRange = $ws.Range("A4:H192")
$sortKey = $ws.Range("B4:B192")
$sortRange.Sort($sortKey, 2)

# --- 3) Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 19:46"
